# Updated symbol list on Wed Jan 18 14:35:31 UTC 2023 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for the crypto rows
# that moved since the last scrape. Values are stored as plain text (the
# sheet already uses text cells for these numeric-looking strings), so we
# force the "@" (Text) number format before assigning each value -- this
# prevents Excel from re-interpreting e.g. "0.003000" as the number 0.003
# or "0.19%" as a native percentage, which would lose the exact formatting
# used by the source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "302.44"
$ws.Range("E2").Value = "0.19%"
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "32.52"
$ws.Range("E3").Value = "0.88%"
$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "5.046"
$ws.Range("E4").Value = "-1.42%"
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07718"
$ws.Range("E5").Value = "-2.16%"
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "2.081"
$ws.Range("E6").Value = "-8.39%"
$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "7.880"
$ws.Range("E7").Value = "0.81%"
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "3.795"
$ws.Range("E8").Value = "-0.24%"
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9258"
$ws.Range("E9").Value = "-0.65%"
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1759"
$ws.Range("E10").Value = "-1.06%"
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07965"
$ws.Range("E11").Value = "3.41%"
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08492"
$ws.Range("E12").Value = "-4.03%"
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03048"
$ws.Range("E13").Value = "-1.06%"
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1001"
$ws.Range("E14").Value = "-0.54%"
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001522"
$ws.Range("E15").Value = "0.60%"
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005714"
$ws.Range("E16").Value = "-4.49%"
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "3.466"
$ws.Range("E18").Value = "0.06%"
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = "2.157"
$ws.Range("E19").Value = "-4.17%"
$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3358"
$ws.Range("E20").Value = "2.65%"
$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1300"
$ws.Range("E21").Value = "-2.93%"
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "4.380"
$ws.Range("E22").Value = "2.26%"
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1975"
$ws.Range("E23").Value = "9.13%"
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04527"
$ws.Range("E24").Value = "-1.35%"
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001239"
$ws.Range("E25").Value = "-0.72%"
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004147"
$ws.Range("E26").Value = "-8.28%"
$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001249"
$ws.Range("E27").Value = "-0.29%"
$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01722"
$ws.Range("E39").Value = "-3.69%"
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04704"
$ws.Range("E40").Value = "-1.20%"
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007472"
$ws.Range("E41").Value = "3.04%"
$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1362"
$ws.Range("E42").Value = "-1.48%"
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002234"
$ws.Range("E43").Value = "5.05%"
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01059"
$ws.Range("E44").Value = "5.03%"
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006208"
$ws.Range("E45").Value = "-2.19%"
$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("E46").Value = "-0.26%"
$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6279"
$ws.Range("E47").Value = "-14.37%"
$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003000"
$ws.Range("E48").Value = "-6.42%"
$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("E49").Value = "-0.26%"
$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("E50").Value = "-0.26%"
